$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "36.634.74"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.02%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.053.98"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.07%  "

$ws.Range("E4").Value = "  -0.11%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "246.80"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.26%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.662"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.59%  "

$ws.Range("E7").Value = "  +0.00%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "54.80"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -6.42%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "60.71"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.88%  "

$ws.Range("E10").Value = "  -2.76%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0754"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.62%  "

$ws.Range("E12").Value = "  -3.04%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.972"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +9.94%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "14.80"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.22%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.353.18"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.13%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.48"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.81%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.068.19"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.11%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "36.554.21"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.13%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "17.14"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -5.61%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "72.04"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.35%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0861"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.16%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "238.37"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.13%  "

$ws.Range("E23").Value = "  -3.46%  "

$ws.Range("E24").Value = "  +0.08%  "

$ws.Range("E25").Value = "  -2.24%  "

$ws.Range("E26").Value = "  +5.55%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "165.54"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.51%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.23"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -9.29%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "20.12"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.03%  "

$ws.Range("E30").Value = "  -1.78%  "

$ws.Range("E31").Value = "  +8.01%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.07"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -7.67%  "

$ws.Range("E33").Value = "  -3.77%  "

$ws.Range("E34").Value = "  -3.57%  "

$ws.Range("B35").Value = "Kaspa"
$ws.Range("C35").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0873"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.21%  "

$ws.Range("B36").Value = "BinanceUSD"
$ws.Range("C36").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.00"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.11%  "

$ws.Range("B37").Value = "WEMIXToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.83"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.33%  "

$ws.Range("B38").Value = "LidoDAOToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.25"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.60%  "

$ws.Range("B39").Value = "TrustWalletToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.25"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -5.81%  "

$ws.Range("B40").Value = "THORChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.04"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.56%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.89"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -5.10%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0215"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.78%  "

$ws.Range("E43").Value = "  -4.76%  "

$ws.Range("E44").Value = "  -2.75%  "

$ws.Range("E45").Value = "  -4.54%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.417.60"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +8.97%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "15.98"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.80%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.50"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +10.35%  "

$ws.Range("E49").Value = "  +1.57%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.27"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.73%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "45.79"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.91%  "
